# Add the new "Gold Price" worksheet after the existing "Currency" sheet.
$wb = $excel.ActiveWorkbook
$currency = $wb.Worksheets.Item("Currency")
$ws = $wb.Worksheets.Add($null, $currency)
$ws.Name = "Gold Price"

# Header row.
$ws.Range("A1").Value = "Date"
$ws.Range("B1").Value = "Price"
$ws.Range("C1").Value = "Change"
$ws.Range("D1").Value = "%Change"

# Match the alignment styling used on the "Currency" sheet:
#   column A -> vertical-center only, columns B:D -> right + vertical-center.
$ws.Range("A1").VerticalAlignment = -4108
$ws.Range("B1:D1").VerticalAlignment = -4108
$ws.Range("B1:D1").HorizontalAlignment = -4152

# Placeholder data rows (2-10), vertical-center aligned like the header.
$ws.Range("A2:D10").VerticalAlignment = -4108

# Column B sized to fit the "Price" values.
$ws.Range("B:B").ColumnWidth = 8.25

# Put the cursor where the author left it on the new sheet.
$ws.Range("E7").Select() | Out-Null
